# Refresh the "Weather Data" sheet with a new reading for Sofia
# (the sheet previously held a reading for Pleven).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Weather Data")

# Helper: write $text into $cellAddr as literal text, even when $text looks
# like a number/percentage to Excel's input parser (e.g. "42.6975", "42%").
# A direct `$ws.Range(...).Value = $text` assignment would silently be
# reinterpreted as a Number for such strings, changing both the stored type
# and (via NumberFormat tricks) the cell's style - neither of which happens
# in the real edit we're reproducing. Routing the text through a formula
# cell ("="..."" always yields a text result) and pasting only the value
# back keeps the literal text *and* the cell's original formatting/style
# intact.
function Set-LiteralText($cellAddr, $text) {
    $escaped = $text -replace '"', '""'
    $scratch = $ws.Range("Z1")
    $scratch.Formula = '="' + $escaped + '"'
    $scratch.Copy() | Out-Null
    $ws.Range($cellAddr).PasteSpecial(-4163) | Out-Null
    $scratch.Clear() | Out-Null
    $excel.CutCopyMode = $false
}

$ws.Range("B2").Value = "Sofia"
Set-LiteralText "B3" "42.6975"
Set-LiteralText "B4" "23.3242"
$ws.Range("B5").Value = "Clear"
$ws.Range("B6").Value = "clear sky"
$ws.Range("B7").Value = "7.83 °C (Feels like 7.83 °C)"
$ws.Range("B8").Value = "7.83 °C to 7.83 °C"
$ws.Range("B9").Value = "1032 hPa"
Set-LiteralText "B10" "42%"
$ws.Range("B11").Value = "1.03 m/s at 0°"
Set-LiteralText "B12" "0%"
$ws.Range("B14").Value = "06:07"
$ws.Range("B15").Value = "16:13"
